$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.311.25'
$ws.Cells.Item(2, 5).Value = '  +0.19%  '
$ws.Cells.Item(3, 4).Value = '1.870.68'
$ws.Cells.Item(3, 5).Value = '  +0.35%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).Value = '235.23'
$ws.Cells.Item(5, 5).Value = '  -0.80%  '
$ws.Cells.Item(6, 4).Value = '1.001'
$ws.Cells.Item(6, 5).Value = '  +0.02%  '
$ws.Cells.Item(7, 4).Value = '0.4697'
$ws.Cells.Item(7, 5).Value = '  +0.26%  '
$ws.Cells.Item(8, 4).Value = '0.2865'
$ws.Cells.Item(8, 5).Value = '  +0.01%  '
$ws.Cells.Item(9, 4).Value = '0.06589'
$ws.Cells.Item(9, 5).Value = '  +0.48%  '
$ws.Cells.Item(10, 4).Value = '21.81'
$ws.Cells.Item(10, 5).Value = '  -1.66%  '
$ws.Cells.Item(11, 4).Value = '0.07994'
$ws.Cells.Item(11, 5).Value = '  +1.15%  '
$ws.Cells.Item(12, 4).Value = '96.99'
$ws.Cells.Item(12, 5).Value = '  -1.27%  '
$ws.Cells.Item(13, 4).Value = '1.871.85'
$ws.Cells.Item(13, 5).Value = '  +0.33%  '
$ws.Cells.Item(14, 4).Value = '0.6896'
$ws.Cells.Item(14, 5).Value = '  +1.20%  '
$ws.Cells.Item(15, 4).Value = '5.119'
$ws.Cells.Item(15, 5).Value = '  -1.25%  '
$ws.Cells.Item(16, 4).Value = '269.40'
$ws.Cells.Item(16, 5).Value = '  -3.07%  '
$ws.Cells.Item(17, 4).Value = '30.337.61'
$ws.Cells.Item(17, 5).Value = '  +0.28%  '
$ws.Cells.Item(18, 5).Value = '  +3.88%  '
$ws.Cells.Item(19, 4).Value = '0.000007655'
$ws.Cells.Item(19, 5).Value = '  +4.22%  '
$ws.Cells.Item(20, 4).Value = '1.000'
$ws.Cells.Item(20, 5).Value = '  -0.01%  '
$ws.Cells.Item(21, 2).Value = 'BinanceUSD'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(21, 4).Value = '1.0000'
$ws.Cells.Item(21, 5).Value = '  -0.12%  '
$ws.Cells.Item(22, 2).Value = 'Uniswap'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(22, 4).Value = '5.269'
$ws.Cells.Item(22, 5).Value = '  -1.79%  '
$ws.Cells.Item(23, 2).Value = 'Chainlink'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(23, 4).Value = '6.223'
$ws.Cells.Item(23, 5).Value = '  +0.41%  '
$ws.Cells.Item(24, 2).Value = 'Cosmos'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(24, 4).Value = '9.394'
$ws.Cells.Item(24, 5).Value = '  +1.73%  '
$ws.Cells.Item(25, 2).Value = 'Monero'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(25, 4).Value = '167.59'
$ws.Cells.Item(25, 5).Value = '  -0.24%  '
$ws.Cells.Item(26, 2).Value = 'EthereumClassic'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(26, 4).Value = '18.89'
$ws.Cells.Item(26, 5).Value = '  -0.97%  '
$ws.Cells.Item(27, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(27, 4).Value = '1.950'
$ws.Cells.Item(27, 5).Value = '  -0.07%  '
$ws.Cells.Item(28, 2).Value = 'Toncoin'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(28, 4).Value = '1.368'
$ws.Cells.Item(28, 5).Value = '  -1.18%  '
$ws.Cells.Item(29, 2).Value = 'Stellar'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(29, 4).Value = '0.09876'
$ws.Cells.Item(29, 5).Value = '  +0.26%  '
$ws.Cells.Item(30, 2).Value = 'Filecoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(30, 4).Value = '4.356'
$ws.Cells.Item(30, 5).Value = '  -0.62%  '
$ws.Cells.Item(31, 2).Value = 'PancakeSwap'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(31, 4).Value = '1.457'
$ws.Cells.Item(31, 5).Value = '  -1.66%  '
$ws.Cells.Item(32, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(32, 4).Value = '4.069'
$ws.Cells.Item(32, 5).Value = '  -0.07%  '
$ws.Cells.Item(33, 2).Value = 'Hedera'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(33, 4).Value = '0.04719'
$ws.Cells.Item(33, 5).Value = '  -0.61%  '
$ws.Cells.Item(34, 2).Value = 'ARBITRUM'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(34, 4).Value = '1.138'
$ws.Cells.Item(34, 5).Value = '  -0.03%  '
$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(35, 4).Value = '0.7027'
$ws.Cells.Item(35, 5).Value = '  -0.26%  '
$ws.Cells.Item(36, 2).Value = 'HuobiToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(36, 4).Value = '2.740'
$ws.Cells.Item(36, 5).Value = '  +1.18%  '
$ws.Cells.Item(37, 2).Value = 'VeChain'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(37, 4).Value = '0.01883'
$ws.Cells.Item(37, 5).Value = '  +0.14%  '
$ws.Cells.Item(38, 2).Value = 'MXToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(38, 4).Value = '2.824'
$ws.Cells.Item(38, 5).Value = '  +7.39%  '
$ws.Cells.Item(39, 2).Value = 'FraxShare'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(39, 4).Value = '6.264'
$ws.Cells.Item(39, 5).Value = '  -0.39%  '
$ws.Cells.Item(40, 2).Value = 'Aave'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(40, 4).Value = '72.23'
$ws.Cells.Item(40, 5).Value = '  -4.68%  '
$ws.Cells.Item(41, 2).Value = 'RenderToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(41, 4).Value = '1.963'
$ws.Cells.Item(41, 5).Value = '  +0.30%  '
$ws.Cells.Item(42, 2).Value = 'TheSandbox'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(42, 4).Value = '0.4181'
$ws.Cells.Item(42, 5).Value = '  +0.14%  '
$ws.Cells.Item(43, 4).Value = '0.8431'
$ws.Cells.Item(43, 5).Value = '  -1.45%  '
$ws.Cells.Item(44, 2).Value = 'PaxDollar'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(44, 4).Value = '1.000'
$ws.Cells.Item(44, 5).Value = '  +0.03%  '
$ws.Cells.Item(45, 2).Value = 'Quant'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(45, 4).Value = '103.10'
$ws.Cells.Item(45, 5).Value = '  -0.30%  '
$ws.Cells.Item(46, 2).Value = 'Aptos'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(46, 4).Value = '7.094'
$ws.Cells.Item(46, 5).Value = '  -1.85%  '
$ws.Cells.Item(47, 2).Value = 'EnergySwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(47, 4).Value = '9.132'
$ws.Cells.Item(47, 5).Value = '  -1.11%  '
$ws.Cells.Item(48, 2).Value = 'Maker'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(48, 4).Value = '918.20'
$ws.Cells.Item(48, 5).Value = '  -3.10%  '
$ws.Cells.Item(49, 2).Value = 'Elrond'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Cells.Item(49, 4).Value = '34.51'
$ws.Cells.Item(49, 5).Value = '  +0.74%  '
$ws.Cells.Item(50, 2).Value = 'Cronos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(50, 4).Value = '0.05699'
$ws.Cells.Item(50, 5).Value = '  +0.93%  '
$ws.Cells.Item(51, 2).Value = 'Algorand'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(51, 4).Value = '0.1120'
$ws.Cells.Item(51, 5).Value = '  -0.20%  '
